$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.932514619883041
$ws.Range("B3").Value = 0.06748538011695907
$ws.Range("B4").Value = 0.9190577399929153
$ws.Range("B5").Value = 0.9586776859504132
$ws.Range("B6").Value = 0.9386586959226578
$ws.Range("B7").Value = 0.952203097417026
